$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer")

$ws.Range("A3").Value = "npeart@mail.com"
$ws.Range("B3").Value = "Neil"
$ws.Range("C3").Value = "Peart"
$ws.Range("D3").Value = "npeart"
$ws.Range("E3").Value = "npeart123"

$ws.Range("C7").Select()
